# Daily auto push: insert one new data row for 2026/01/17 13:00 at row 655,
# pushing the existing rows 655:696 down to 656:697.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(655).Insert()

# Force column A to be stored as plain text so the "yyyy/mm/dd" string isn't
# auto-converted into a date serial number, then drop back to the default
# "Normal" style so no stray number-format style sticks to the cell.
$ws.Range("A655").NumberFormat = "@"
$ws.Range("A655").Value = "2026/01/17"
$ws.Range("A655").Style = "Normal"
$ws.Range("B655").Value = "土"
$ws.Range("C655").Value = 13
$ws.Range("D655").Value = 201
